$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new displayed text, taken from the commit diff (cryptos.xlsx row refresh).
$updates = @(
    @{ Cell = "D2"; Value = "26.351.14" }
    @{ Cell = "E2"; Value = "  +0.38%  " }
    @{ Cell = "D3"; Value = "1.694.86" }
    @{ Cell = "E3"; Value = "  +0.93%  " }
    @{ Cell = "E4"; Value = "  +0.15%  " }
    @{ Cell = "D5"; Value = "218.01" }
    @{ Cell = "E5"; Value = "  -0.25%  " }
    @{ Cell = "D6"; Value = "0.5416" }
    @{ Cell = "E6"; Value = "  +2.66%  " }
    @{ Cell = "E7"; Value = "  +0.12%  " }
    @{ Cell = "E8"; Value = "  +1.09%  " }
    @{ Cell = "D9"; Value = "0.06456" }
    @{ Cell = "E9"; Value = "  -0.50%  " }
    @{ Cell = "E10"; Value = "  -1.33%  " }
    @{ Cell = "D11"; Value = "0.07661" }
    @{ Cell = "E11"; Value = "  +1.78%  " }
    @{ Cell = "D12"; Value = "1.728.08" }
    @{ Cell = "E12"; Value = "  +2.74%  " }
    @{ Cell = "D13"; Value = "4.539" }
    @{ Cell = "E13"; Value = "  +0.25%  " }
    @{ Cell = "D14"; Value = "0.5816" }
    @{ Cell = "E14"; Value = "  +0.10%  " }
    @{ Cell = "D15"; Value = "0.000008413" }
    @{ Cell = "E15"; Value = "  -1.21%  " }
    @{ Cell = "D16"; Value = "67.13" }
    @{ Cell = "E16"; Value = "  +3.85%  " }
    @{ Cell = "D17"; Value = "26.400.71" }
    @{ Cell = "E17"; Value = "  +0.35%  " }
    @{ Cell = "D18"; Value = "4.922" }
    @{ Cell = "E18"; Value = "  -0.10%  " }
    @{ Cell = "E19"; Value = "  +0.10%  " }
    @{ Cell = "E20"; Value = "  +0.25%  " }
    @{ Cell = "E21"; Value = "  +0.33%  " }
    @{ Cell = "D22"; Value = "6.284" }
    @{ Cell = "E22"; Value = "  +1.35%  " }
    @{ Cell = "E23"; Value = "  +0.10%  " }
    @{ Cell = "D24"; Value = "149.14" }
    @{ Cell = "E24"; Value = "  +2.52%  " }
    @{ Cell = "D25"; Value = "0.1285" }
    @{ Cell = "E25"; Value = "  +3.38%  " }
    @{ Cell = "D26"; Value = "7.857" }
    @{ Cell = "D27"; Value = "15.89" }
    @{ Cell = "E27"; Value = "  +0.47%  " }
    @{ Cell = "D28"; Value = "0.06335" }
    @{ Cell = "E28"; Value = "  -3.50%  " }
    @{ Cell = "D29"; Value = "1.386" }
    @{ Cell = "E29"; Value = "  +2.36%  " }
    @{ Cell = "E30"; Value = "  -0.33%  " }
    @{ Cell = "D31"; Value = "3.615" }
    @{ Cell = "E31"; Value = "  +0.65%  " }
    @{ Cell = "D32"; Value = "3.593" }
    @{ Cell = "E32"; Value = "  -0.22%  " }
    @{ Cell = "D34"; Value = "1.033" }
    @{ Cell = "E34"; Value = "  +0.07%  " }
    @{ Cell = "D35"; Value = "0.6202" }
    @{ Cell = "E35"; Value = "  -0.45%  " }
    @{ Cell = "D36"; Value = "2.416" }
    @{ Cell = "E36"; Value = "  +0.54%  " }
    @{ Cell = "D37"; Value = "2.754" }
    @{ Cell = "E37"; Value = "  +0.78%  " }
    @{ Cell = "D38"; Value = "0.01655" }
    @{ Cell = "E38"; Value = "  +1.94%  " }
    @{ Cell = "D39"; Value = "1.115.92" }
    @{ Cell = "E39"; Value = "  +0.29%  " }
    @{ Cell = "D40"; Value = "6.106" }
    @{ Cell = "E40"; Value = "  -5.37%  " }
    @{ Cell = "E41"; Value = "  +1.30%  " }
    @{ Cell = "E42"; Value = "  +0.02%  " }
    @{ Cell = "D43"; Value = "101.11" }
    @{ Cell = "E43"; Value = "  +0.31%  " }
    @{ Cell = "D44"; Value = "1.846.85" }
    @{ Cell = "E44"; Value = "  +0.97%  " }
    @{ Cell = "E45"; Value = "  +1.22%  " }
    @{ Cell = "D46"; Value = "57.77" }
    @{ Cell = "E46"; Value = "  +1.27%  " }
    @{ Cell = "D47"; Value = "8.187" }
    @{ Cell = "E47"; Value = "  -0.06%  " }
    @{ Cell = "D48"; Value = "1.003" }
    @{ Cell = "E48"; Value = "  -0.35%  " }
    @{ Cell = "D49"; Value = "0.05284" }
    @{ Cell = "E49"; Value = "  +0.22%  " }
    @{ Cell = "B50"; Value = "Mantle" }
    @{ Cell = "C50"; Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt" }
    @{ Cell = "D50"; Value = "0.4303" }
    @{ Cell = "E50"; Value = "  +0.26%  " }
    @{ Cell = "B51"; Value = "Aptos" }
    @{ Cell = "C51"; Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt" }
    @{ Cell = "D51"; Value = "6.101" }
    @{ Cell = "E51"; Value = "  +0.16%  " }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $text = $u.Value
    # Column D holds price strings such as "26.351.14" or "218.01" that must
    # stay plain text (they are not real numbers - some use "." as a thousands
    # separator). Writing them straight to .Value lets Excel coerce anything
    # that parses cleanly (e.g. "218.01") into a floating-point number, which
    # then round-trips with binary noise (218.00999999999999). Force text via
    # NumberFormat "@", then restore the cells original style so no stray
    # formatting is left behind.
    if ($u.Cell.StartsWith("D")) {
        $originalStyle = $cell.Style
        $cell.NumberFormat = "@"
        $cell.Value = $text
        $cell.Style = $originalStyle
    } else {
        $cell.Value = $text
    }
}
